$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 3800
$ws.Range("E2").Value = 361
$ws.Range("F2").Value = 395
$ws.Range("G2").Value = 422
$ws.Range("H2").Value = 420
$ws.Range("I2").Value = 423
$ws.Range("J2").Value = -3
$ws.Range("K2").Value = 4575
$ws.Range("L2").Value = 1389
$ws.Range("M2").Value = 3186
$ws.Range("N2").Value = 3183
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 132
$ws.Range("Q2").Value = -747
$ws.Range("R2").Value = -462
$ws.Range("S2").Value = 479
$ws.Range("T2").Value = 274
$ws.Range("U2").Value = -1021
$ws.Range("V2").Value = 1005
$ws.Range("W2").Value = 9.51
$ws.Range("X2").Value = 11.05
$ws.Range("Y2").Value = 14.34
$ws.Range("Z2").Value = 10.5
$ws.Range("AA2").Value = 43.59
$ws.Range("AB2").Value = 2381.07
$ws.Range("AC2").Value = 827
$ws.Range("AD2").Value = 22.2
$ws.Range("AE2").Value = 6229
$ws.Range("AF2").Value = 2.95
$ws.Range("AG2").Value = 19
$ws.Range("AH2").Value = 0.1
$ws.Range("AJ2").Value = 51120439
$ws.Range("AI2").ClearContents()

# Row 3
$ws.Range("D3").Value = 4445
$ws.Range("E3").Value = 335
$ws.Range("F3").Value = 477
$ws.Range("G3").Value = 483
$ws.Range("H3").Value = 469
$ws.Range("I3").Value = 476
$ws.Range("J3").Value = -7
$ws.Range("K3").Value = 6002
$ws.Range("L3").Value = 2149
$ws.Range("M3").Value = 3853
$ws.Range("N3").Value = 3858
$ws.Range("O3").Value = -4
$ws.Range("P3").Value = 140
$ws.Range("Q3").Value = -186
$ws.Range("R3").Value = -261
$ws.Range("S3").Value = 310
$ws.Range("T3").Value = 147
$ws.Range("U3").Value = -334
$ws.Range("V3").Value = 1433
$ws.Range("W3").Value = 7.54
$ws.Range("X3").Value = 10.54
$ws.Range("Y3").Value = 13.51
$ws.Range("Z3").Value = 8.859999999999999
$ws.Range("AA3").Value = 55.78
$ws.Range("AB3").Value = 2728.4
$ws.Range("AC3").Value = 930
$ws.Range("AD3").Value = 10.38
$ws.Range("AE3").Value = 7549
$ws.Range("AF3").Value = 1.28
$ws.Range("AG3").Value = 94
$ws.Range("AH3").Value = 0.97
$ws.Range("AJ3").Value = 51120439
$ws.Range("AI3").ClearContents()

# Row 4
$ws.Range("D4").Value = 3923
$ws.Range("E4").Value = 205
$ws.Range("F4").Value = 564
$ws.Range("G4").Value = 419
$ws.Range("H4").Value = 416
$ws.Range("I4").Value = 430
$ws.Range("J4").Value = -14
$ws.Range("K4").Value = 6901
$ws.Range("L4").Value = 2520
$ws.Range("M4").Value = 4381
$ws.Range("N4").Value = 4399
$ws.Range("O4").Value = -19
$ws.Range("P4").Value = 145
$ws.Range("Q4").Value = 36
$ws.Range("R4").Value = -157
$ws.Range("S4").Value = 58
$ws.Range("T4").Value = 107
$ws.Range("U4").Value = -71
$ws.Range("V4").Value = 1586
$ws.Range("W4").Value = 5.23
$ws.Range("X4").Value = 10.61
$ws.Range("Y4").Value = 10.41
$ws.Range("Z4").Value = 6.45
$ws.Range("AA4").Value = 57.53
$ws.Range("AB4").Value = 3004.03
$ws.Range("AC4").Value = 841
$ws.Range("AD4").Value = 8.94
$ws.Range("AE4").Value = 8610
$ws.Range("AF4").Value = 0.87
$ws.Range("AG4").Value = 112
$ws.Range("AH4").Value = 1.49
$ws.Range("AJ4").Value = 51120439
$ws.Range("AI4").ClearContents()

# Row 5
$ws.Range("D5").Value = 3332
$ws.Range("E5").Value = -226
$ws.Range("F5").Value = 376
$ws.Range("G5").Value = 140
$ws.Range("H5").Value = 136
$ws.Range("I5").Value = 149
$ws.Range("J5").Value = -13
$ws.Range("K5").Value = 6306
$ws.Range("L5").Value = 2329
$ws.Range("M5").Value = 3977
$ws.Range("N5").Value = 3988
$ws.Range("O5").Value = -11
$ws.Range("P5").Value = 128
$ws.Range("Q5").Value = 168
$ws.Range("R5").Value = -190
$ws.Range("S5").Value = 192
$ws.Range("T5").Value = 125
$ws.Range("U5").Value = 43
$ws.Range("V5").Value = 1720
$ws.Range("W5").Value = -6.79
$ws.Range("X5").Value = 4.08
$ws.Range("Y5").Value = 3.55
$ws.Range("Z5").Value = 2.06
$ws.Range("AA5").Value = 58.56
$ws.Range("AB5").Value = 3071.56
$ws.Range("AC5").Value = 291
$ws.Range("AD5").Value = 15.22
$ws.Range("AE5").Value = 7805
$ws.Range("AF5").Value = 0.57
$ws.Range("AG5").Value = 112
$ws.Range("AH5").Value = 2.53
$ws.Range("AJ5").Value = 51120439
$ws.Range("AI5").ClearContents()

# Row 6
$ws.Range("D6").Value = 2940
$ws.Range("E6").Value = -615
$ws.Range("F6").Value = -148
$ws.Range("G6").Value = -487
$ws.Range("H6").Value = -490
$ws.Range("I6").Value = -444
$ws.Range("K6").Value = 6537
$ws.Range("L6").Value = 2940
$ws.Range("M6").Value = 3596
$ws.Range("N6").Value = 3581
$ws.Range("P6").Value = 134
$ws.Range("Q6").Value = 7
$ws.Range("R6").Value = -361
$ws.Range("S6").Value = 373
$ws.Range("T6").Value = 143
$ws.Range("U6").Value = -136
$ws.Range("V6").Value = 2312
$ws.Range("W6").Value = -20.93
$ws.Range("X6").Value = -16.68
$ws.Range("Y6").Value = -11.73
$ws.Range("Z6").Value = -7.64
$ws.Range("AA6").Value = 81.75
$ws.Range("AB6").Value = 2670.21
$ws.Range("AC6").Value = -869
$ws.Range("AD6").Value = -2.68
$ws.Range("AE6").Value = 7013
$ws.Range("AF6").Value = 0.33
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AJ6").Value = 51120439
$ws.Range("AI6").ClearContents()

# Row 7: clear all data cells except A, B, C
$ws.Range("D7:AI7").ClearContents()

# Row 8: clear all data cells except A, B, C
$ws.Range("D8:AI8").ClearContents()

# Row 9: clear all data cells except A, B, C
$ws.Range("D9:AI9").ClearContents()

Write-Host "Edits applied"